# Fix "Recorded By" (column G) entries: move a leading "System, " tag
# from the front of the comma-separated list to the end, e.g.
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com"   -> "system, backup@backdoor.com, System"
# Entries that don't start with "System, " (including the bare word
# "System" with nothing after it) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$prefix = "System, "

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($null -ne $val -and $val.StartsWith($prefix)) {
        $rest = $val.Substring($prefix.Length)
        $cell.Value = "$rest, System"
    }
}
